# Update the "Förändrad" (Changed) date column (C) from 45203 to 45204
# for all data rows (rows 2 through 78) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 78; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value2 = 45204
    }
}
